# Auto-generated Excel COM-interop script to apply symbol-list price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking values in column D must stay stored as text (matching the
# original inline-string cells), so force each target cell to Text format
# before writing the value - otherwise Excel auto-converts them to numbers.
$priceCells = @(
    "D2"
    "D4"
    "D5"
    "D6"
    "D7"
    "D8"
    "D9"
    "D10"
    "D11"
    "D13"
    "D14"
    "D15"
    "D16"
    "D17"
    "D18"
    "D19"
    "D21"
    "D22"
    "D23"
    "D24"
    "D26"
    "D27"
    "D40"
    "D41"
    "D42"
    "D43"
    "D44"
    "D45"
    "D47"
    "D48"
    "D49"
    "D50"
)
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Updated price values (column D)
$ws.Range("D2").Value2 = "243.36"
$ws.Range("D4").Value2 = "5.404"
$ws.Range("D5").Value2 = "0.05982"
$ws.Range("D6").Value2 = "3.422"
$ws.Range("D7").Value2 = "6.495"
$ws.Range("D8").Value2 = "0.8131"
$ws.Range("D9").Value2 = "0.9309"
$ws.Range("D10").Value2 = "0.1434"
$ws.Range("D11").Value2 = "0.07450"
$ws.Range("D13").Value2 = "0.03038"
$ws.Range("D14").Value2 = "0.09353"
$ws.Range("D15").Value2 = "3.851"
$ws.Range("D16").Value2 = "0.001589"
$ws.Range("D17").Value2 = "0.04695"
$ws.Range("D18").Value2 = "0.0005902"
$ws.Range("D19").Value2 = "0.005880"
$ws.Range("D21").Value2 = "0.004877"
$ws.Range("D22").Value2 = "0.00008002"
$ws.Range("D23").Value2 = "3.571"
$ws.Range("D24").Value2 = "2.175"
$ws.Range("D26").Value2 = "0.1302"
$ws.Range("D27").Value2 = "0.0002340"
$ws.Range("D40").Value2 = "0.03957"
$ws.Range("D41").Value2 = "0.006382"
$ws.Range("D42").Value2 = "0.004001"
$ws.Range("D43").Value2 = "0.1075"
$ws.Range("D44").Value2 = "0.008569"
$ws.Range("D45").Value2 = "0.00005177"
$ws.Range("D47").Value2 = "0.6702"
$ws.Range("D48").Value2 = "0.002267"
$ws.Range("D49").Value2 = "0.00002101"
$ws.Range("D50").Value2 = "0.0002001"

# Updated volume/label text values (column E)
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E22").Value = "21NitroExNTXWorstin24h"
